$d = $word.ActiveDocument

# Paragraph 2: "PHI (Personal Health Information (Ref-DJ49F2))" -> "PHI (Ref-f104270))"
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("PHI (Personal Health Information (Ref-DJ49F2))", $true, $false, $false, $false, $false, $true, 1, $false, "PHI (Ref-f104270))", 2)

# Paragraph 3: two distinct reference tags both become "Ref-f385834"
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Ref-DJ49F2", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-f385834", 2)
$p3b = $d.Paragraphs(3).Range
$p3b.Find.Execute("Ref-G7H21K", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-f385834", 2)

# Paragraph 4: "Ref-J7X2BZ" -> "Lee 208"
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Ref-J7X2BZ", $true, $false, $false, $false, $false, $true, 1, $false, "Lee 208", 2)

# Paragraph 5: "Ref-DJ49F2" -> "Nguyen 59" + en-dash + "60"
$p5 = $d.Paragraphs(5).Range
$enDash = [char]0x2013
$p5.Find.Execute("Ref-DJ49F2", $true, $false, $false, $false, $false, $true, 1, $false, "Nguyen 59${enDash}60", 2)

# Paragraph 7: "Ref-A1B2C3" -> "Ref-f719863"
$p7 = $d.Paragraphs(7).Range
$p7.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false, $true, 1, $false, "Ref-f719863", 2)
